$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two existing trailing rows (B and D values were refreshed
#     upstream) ---
$ws.Range("B325").Value = 7710432000000
$ws.Range("D325").Value = 255905476269.4988

$ws.Range("B326").Value = 7822810000000
$ws.Range("D326").Value = 255814584695.8797

# --- Append three new monthly data rows (327-329), cloning the
#     formatting (date style incl. number format/border/alignment) of the
#     last existing data row so the new date cells look the same ---
$ws.Range("A326:D326").Copy()
$ws.Range("A327:D329").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A327").Value = 44986
$ws.Range("B327").Value = 7965088000000
$ws.Range("C327").Value = 0.03252032520325204
$ws.Range("D327").Value = 259027252032.5204

$ws.Range("A328").Value = 45017
$ws.Range("B328").Value = 8069151000000
$ws.Range("C328").Value = 0.03241491085899514
$ws.Range("D328").Value = 261560810372.7715

$ws.Range("A329").Value = 45047
$ws.Range("B329").Value = 8140535000000
$ws.Range("C329").Value = 0.03241491085899514
$ws.Range("D329").Value = 263874716369.53
